# The workbook has two sheets; the data-collection results live on the
# first (active) sheet - "八位序列号收集收集结果yd5", physically sheet1.xml.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 35 -----------------------------------------------------
$ws.Range("A35").Value = "......"

$ws.Range("B35").NumberFormat = "yyyy/m/d h:mm:ss;@"
$ws.Range("B35").Value = 45909.007025463

$ws.Range("C35").Value = "fd81cb3d"

# D holds an (8-9 digit) QQ number that must stay textual, like the
# existing D26:D34 cells, so force the cell to text before typing it.
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "908289638"

# --- New row 36 -----------------------------------------------------
$ws.Range("A36").Value = "YL."

$ws.Range("B36").NumberFormat = "yyyy/m/d h:mm:ss;@"
$ws.Range("B36").Value = 45909.3519560185

$ws.Range("C36").Value = "a4c3e725"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "918357021"
